$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.527.03"

$ws.Cells.Item(3, 4).Value = "2.139.66"
$ws.Cells.Item(3, 5).Value = "  +1.84%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.008"
$ws.Cells.Item(4, 5).Value = "  +0.55%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "352.24"
$ws.Cells.Item(5, 5).Value = "  +5.21%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.5257"
$ws.Cells.Item(7, 5).Value = "  +0.64%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.4566"
$ws.Cells.Item(8, 5).Value = "  -0.11%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "53.66"
$ws.Cells.Item(9, 5).Value = "  -3.72%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.09184"
$ws.Cells.Item(10, 5).Value = "  +3.02%  "

$ws.Cells.Item(11, 5).Value = "  +1.35%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "25.57"
$ws.Cells.Item(12, 5).Value = "  +5.74%  "

$ws.Cells.Item(13, 4).Value = "2.125.00"
$ws.Cells.Item(13, 5).Value = "  +1.33%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.901"
$ws.Cells.Item(14, 5).Value = "  +1.36%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "8.175"
$ws.Cells.Item(15, 5).Value = "  +1.53%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "101.59"
$ws.Cells.Item(16, 5).Value = "  +4.63%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.00001170"
$ws.Cells.Item(17, 5).Value = "  +2.07%  "

$ws.Cells.Item(18, 5).Value = "  +0.35%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06724"
$ws.Cells.Item(19, 5).Value = "  +1.38%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "20.58"
$ws.Cells.Item(20, 5).Value = "  +7.21%  "

$ws.Cells.Item(21, 5).Value = "  +0.44%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.379"
$ws.Cells.Item(22, 5).Value = "  +1.28%  "

$ws.Cells.Item(23, 4).Value = "30.609.26"
$ws.Cells.Item(23, 5).Value = "  +0.44%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "12.86"
$ws.Cells.Item(24, 5).Value = "  +4.21%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.376"
$ws.Cells.Item(25, 5).Value = "  +0.70%  "

$ws.Cells.Item(26, 4).Value = "2.394.14"
$ws.Cells.Item(26, 5).Value = "  +2.19%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "22.51"
$ws.Cells.Item(27, 5).Value = "  +1.40%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.606"
$ws.Cells.Item(28, 5).Value = "  +3.67%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "165.11"
$ws.Cells.Item(29, 5).Value = "  +1.75%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "135.97"
$ws.Cells.Item(30, 5).Value = "  +1.99%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.224"
$ws.Cells.Item(31, 5).Value = "  +1.50%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.721"
$ws.Cells.Item(32, 5).Value = "  +3.98%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.1083"
$ws.Cells.Item(33, 5).Value = "  +1.40%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "6.409"
$ws.Cells.Item(34, 5).Value = "  +0.29%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.038"
$ws.Cells.Item(35, 5).Value = "  +2.70%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "6.169"
$ws.Cells.Item(36, 5).Value = "  +5.03%  "

$ws.Cells.Item(37, 5).Value = "  +1.55%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02647"
$ws.Cells.Item(38, 5).Value = "  +2.86%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.06996"
$ws.Cells.Item(39, 5).Value = "  +2.03%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.2359"
$ws.Cells.Item(40, 5).Value = "  +1.67%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "12.75"
$ws.Cells.Item(41, 5).Value = "  +0.70%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.6998"
$ws.Cells.Item(42, 5).Value = "  +1.79%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.277"
$ws.Cells.Item(43, 5).Value = "  +2.37%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "14.83"
$ws.Cells.Item(44, 5).Value = "  +5.92%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.354"
$ws.Cells.Item(45, 5).Value = "  +1.30%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.6532"
$ws.Cells.Item(46, 5).Value = "  +2.37%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.00000000377"
$ws.Cells.Item(47, 5).Value = "  +11.19%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.724"
$ws.Cells.Item(48, 5).Value = "  +1.82%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.252"
$ws.Cells.Item(49, 5).Value = "  +0.18%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "83.91"
$ws.Cells.Item(50, 5).Value = "  +0.93%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.07300"
$ws.Cells.Item(51, 5).Value = "  +2.47%  "
